$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Current": insert new rows describing the ACS723 bias, and shift the
# existing "potential divider" block down to make room.
# ---------------------------------------------------------------------------
$wsCurrent = $wb.Worksheets.Item("Current")

# Insert 4 new rows starting at row 16 (pushes old rows 16-19 to 20-23).
$wsCurrent.Range("A16:A19").EntireRow.Insert()

$wsCurrent.Cells.Item(16, 1).Value = "ACS723 bias"
$wsCurrent.Cells.Item(16, 2).Value = 0.5
$wsCurrent.Cells.Item(16, 3).Value = "V"

$wsCurrent.Cells.Item(17, 1).Value = "bias at input to Ganymede"
$wsCurrent.Cells.Item(17, 2).Formula = "=B16*B9"
$wsCurrent.Cells.Item(17, 3).Value = "V"

$wsCurrent.Cells.Item(18, 1).Value = "equivalent bias current"
$wsCurrent.Cells.Item(18, 2).Formula = "=B17/B10"
$wsCurrent.Cells.Item(18, 3).Value = "A"

# ---------------------------------------------------------------------------
# Sheet "comparator thresholds": update the current trip threshold to 27A and
# add a row that includes the ACS723 bias in the comparator trip current.
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("comparator thresholds")

# Update the trip current value itself.
$wsComp.Cells.Item(10, 2).Value = 27

# Insert a new row after "trip current" (row 10) for the bias-inclusive trip.
$wsComp.Rows.Item(11).Insert()

$wsComp.Cells.Item(11, 1).Value = "trip including bias"
$wsComp.Cells.Item(11, 2).Formula = "=B10+Current!B18"
$wsComp.Cells.Item(11, 3).Value = "A apparent"

# The formerly-row-11 "comparator voltage" (now row 12) should reference the
# new "trip including bias" row instead of the raw trip current.
$wsComp.Cells.Item(12, 2).Formula = "=B11*Current!B10"

# The "DAC setting" row (now row 13) picks up the same integer display style
# used by the other DAC setting rows on this sheet.
$wsComp.Cells.Item(13, 2).NumberFormat = "0"

# Make "comparator thresholds" the active sheet/tab, matching the workbook.
$wsComp.Activate()
$wsComp.Range("B11").Select()
